$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scrub patient-identifying notes (HPV / medication / DOB strings) from column I,
# replacing them (and the now-redundant "HPV patient" label) with a single blank space.
for ($r = 3; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = " "
}

# Touch A20's formatting so it collapses onto the equivalent, already-existing style
# (font size re-applied as a no-op) instead of keeping its own duplicate style entry.
$ws.Range("A20").Font.Size = 14

# Move the active selection to I11, reflecting where the user left off editing.
$ws.Range("I11").Select() | Out-Null
